$wb = $excel.ActiveWorkbook

$wsMovies = $wb.Worksheets.Item("movies")
$wsTv = $wb.Worksheets.Item("live-action-tv-series")
$wsAnim = $wb.Worksheets.Item("animated-tv-series")

# --- live-action-tv-series: drop the boolean "mcu" column (F), keep the
#     richer "continuity" column (was G) by deleting column F so G shifts
#     left into F.
$wsTv.Columns.Item(6).Delete()

# --- Reclassify the Netflix-era Marvel shows from "mcu" to "mcu-netflix"
#     (Agent Carter, Daredevil, Jessica Jones, Luke Cage, Iron Fist,
#     The Defenders, The Punisher).
$netflixRows = @(4, 5, 6, 7, 9, 10, 13)
foreach ($r in $netflixRows) {
    $wsTv.Range("F" + $r).Value = "mcu-netflix"
}

# --- Selections / active-cell bookkeeping -------------------------------
# movies: selection moves from E69 to C69, and the sheet is no longer the
# tab that's active/selected.
$wsMovies.Range("C69").Select()

# live-action-tv-series: selection moves from A25 to C7.
$wsTv.Range("C7").Select()

# animated-tv-series: becomes the active/selected tab; its scroll position
# (topLeftCell) moves from A10 to A11, selection itself is unchanged.
$wsAnim.Range("A10").Select()
$excel.ActiveWindow.ScrollRow = 11
